$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting so values like
# "1.002" or "29.286.86" are not reinterpreted as numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.286.86'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").Value = '1.858.48'
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = '0.7012'
$ws.Range("E5").Value = '  +1.46%  '
$ws.Range("D6").Value = '238.11'
$ws.Range("E6").Value = '  +0.39%  '
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").Value = '0.07870'
$ws.Range("E8").Value = '  +2.08%  '
$ws.Range("D9").Value = '0.3028'
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").Value = '24.45'
$ws.Range("E10").Value = '  +5.99%  '
$ws.Range("D11").Value = '0.08163'
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("D12").Value = '1.876.87'
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("D13").Value = '5.209'
$ws.Range("E13").Value = '  +1.13%  '
$ws.Range("D14").Value = '0.7064'
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("D15").Value = '89.49'
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").Value = '29.337.78'
$ws.Range("E16").Value = '  +0.63%  '
$ws.Range("D17").Value = '5.805'
$ws.Range("E17").Value = '  +1.63%  '
$ws.Range("D18").Value = '0.000007820'
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("D19").Value = '13.20'
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").Value = '237.07'
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("D21").Value = '2.123.60'
$ws.Range("E21").Value = '  +1.16%  '
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = '7.569'
$ws.Range("E24").Value = '  +1.78%  '
$ws.Range("D25").Value = '162.32'
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("D26").Value = '8.894'
$ws.Range("E26").Value = '  -1.09%  '
$ws.Range("D27").Value = '0.1415'
$ws.Range("E27").Value = '  -0.97%  '
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("D29").Value = '1.904'
$ws.Range("E29").Value = '  -2.32%  '
$ws.Range("D30").Value = '1.397'
$ws.Range("E30").Value = '  -0.62%  '
$ws.Range("D31").Value = '1.480'
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").Value = '4.295'
$ws.Range("E32").Value = '  -3.81%  '
$ws.Range("D33").Value = '4.030'
$ws.Range("E33").Value = '  +1.02%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  +0.97%  '
$ws.Range("D36").Value = '0.7078'
$ws.Range("E36").Value = '  +0.30%  '
$ws.Range("D37").Value = '0.9997'
$ws.Range("E37").Value = '  +0.19%  '
$ws.Range("D38").Value = '2.680'
$ws.Range("E38").Value = '  +1.07%  '
$ws.Range("D39").Value = '0.01849'
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("D40").Value = '2.684'
$ws.Range("E40").Value = '  -1.11%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '0.9214'
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.139.61'
$ws.Range("E42").Value = '  +2.60%  '
$ws.Range("D43").Value = '5.954'
$ws.Range("E43").Value = '  +1.51%  '
$ws.Range("D44").Value = '0.4240'
$ws.Range("E44").Value = '  -0.49%  '
$ws.Range("D45").Value = '70.31'
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("D47").Value = '102.88'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("D48").Value = '0.5319'
$ws.Range("E48").Value = '  -2.53%  '
$ws.Range("D49").Value = '1.744'
$ws.Range("E49").Value = '  -2.41%  '
$ws.Range("D50").Value = '9.183'
$ws.Range("E50").Value = '  +0.60%  '
$ws.Range("D51").Value = '6.979'
$ws.Range("E51").Value = '  +0.46%  '
